$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 103 ---
$ws.Cells.Item(103, 1).Value = 45462.2916666667
$ws.Cells.Item(103, 1).NumberFormat = "yyyy-mm-dd hh:mm:ss"
$ws.Cells.Item(103, 1).Font.Name = "Calibri"
$ws.Cells.Item(103, 1).Font.Color = 0

$ws.Cells.Item(103, 2).Value = 0
$ws.Cells.Item(103, 3).Value = 2
$ws.Cells.Item(103, 4).Value = 2
$ws.Cells.Item(103, 5).Value = 2
$ws.Cells.Item(103, 6).Value = 2

$ws.Cells.Item(103, 7).Formula = '="2"'
$ws.Cells.Item(103, 7).Copy()
$ws.Cells.Item(103, 7).PasteSpecial(-4163)

$ws.Cells.Item(103, 8).Value = "KK.MI"

# --- Row 104 ---
$ws.Cells.Item(104, 1).Value = 45463.4677199074
$ws.Cells.Item(104, 1).NumberFormat = "yyyy-mm-dd hh:mm:ss"
$ws.Cells.Item(104, 1).Font.Name = "Calibri"
$ws.Cells.Item(104, 1).Font.Color = 0

$ws.Cells.Item(104, 2).Value = 1800
$ws.Cells.Item(104, 3).Value = 1.96000003814697
$ws.Cells.Item(104, 4).Value = 1.96000003814697
$ws.Cells.Item(104, 5).Value = 1.96000003814697
$ws.Cells.Item(104, 6).Value = 1.96000003814697

$ws.Cells.Item(104, 7).Formula = '="1.96000003814697"'
$ws.Cells.Item(104, 7).Copy()
$ws.Cells.Item(104, 7).PasteSpecial(-4163)

$ws.Cells.Item(104, 8).Value = "KK.MI"
